# "committing 1 player mode modifications"
#
# The title page (blank bidi paragraph + 5 RTL name/CS-IS/Sec1-Sec2 lines)
# is removed. The surviving "Chess game" heading paragraph becomes the new
# first paragraph of the document: it keeps its centered pPr (sz/szCs 40)
# and its run's sz/szCs 72 formatting, loses its now-irrelevant
# <w:lastRenderedPageBreak/> hint, and gains the "_GoBack" bookmark that
# used to sit later in the document (right before the "Movement-Generation"
# paragraph's first run "Move").

$d = $word.ActiveDocument

# 1) Relocate the "_GoBack" bookmark: delete it from its old spot...
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ...and re-add it collapsed at the very start of the "Chess game" paragraph
# (paragraph 7: blank bidi para + 5 name paragraphs + "Chess game"). Doing
# this before the paragraphs above it are deleted keeps the insertion point
# unambiguous (it sits right before real run content, just like the
# original bookmark did), so both bookmarkStart/bookmarkEnd land inside
# this paragraph instead of drifting into a neighboring one.
$chessPara = $d.Paragraphs.Item(7)
$bmRange = $d.Range($chessPara.Range.Start, $chessPara.Range.Start)
$d.Bookmarks.Add("_GoBack", $bmRange)

# 2) Delete the title-page paragraphs (1 through 6): the empty bidi
# paragraph plus the five RTL name / CS|IS / Sec1|Sec2 lines. Word merges
# the remainder into paragraph 1, which is now the "Chess game" paragraph,
# carrying over its own pPr/rPr untouched.
$firstPara = $d.Paragraphs.Item(1)
$lastPara = $d.Paragraphs.Item(6)
$titleRange = $d.Range($firstPara.Range.Start, $lastPara.Range.End)
$titleRange.Delete()

# 3) Re-stamp the paragraph's text in place. This is a no-op textually
# ("Chess game" -> "Chess game") but it drops the stale
# <w:lastRenderedPageBreak/> rendering hint left over from the old layout,
# while leaving the paragraph/run formatting (jc center; sz/szCs 40 on the
# paragraph mark; sz/szCs 72 on the run) intact.
$chessParaNow = $d.Paragraphs.Item(1)
$chessParaNow.Range.Text = "Chess game"
